$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
$tcs.Name = "Office"
Write-Host ("after set: " + $tcs.Name)
